$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string (e.g. "14.16") need the
# columns NumberFormat forced to Text first, otherwise Excel auto-converts the
# assigned string into a numeric cell instead of leaving it as text.
$textFormatCells = @(
    "D5",
    "D6",
    "D8",
    "D10",
    "D11",
    "D14",
    "D15",
    "D16",
    "D17",
    "D20",
    "D22",
    "D23",
    "D27",
    "D28",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D40",
    "D42",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$newValues = [ordered]@{
    "D2" = "36.504.93"
    "E2" = "  -2.14%  "
    "D3" = "1.992.48"
    "E3" = "  -0.77%  "
    "E4" = "  -0.02%  "
    "D5" = "233.26"
    "E5" = "  -9.59%  "
    "D6" = "0.599"
    "E6" = "  -1.82%  "
    "E7" = "  +0.04%  "
    "D8" = "54.81"
    "E8" = "  -1.80%  "
    "E9" = "  -4.11%  "
    "D10" = "57.83"
    "E10" = "  +2.49%  "
    "D11" = "0.0745"
    "E11" = "  -3.00%  "
    "E12" = "  -3.14%  "
    "B13" = "WrappedliquidstakedEther2.0"
    "C13" = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
    "D13" = "2.290.00"
    "E13" = "  -0.60%  "
    "B14" = "Chainlink"
    "C14" = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
    "D14" = "14.16"
    "E14" = "  -0.34%  "
    "D15" = "20.21"
    "E15" = "  -3.23%  "
    "D16" = "0.756"
    "E16" = "  -5.44%  "
    "D17" = "5.06"
    "E17" = "  -3.27%  "
    "D18" = "2.011.53"
    "E18" = "  +0.48%  "
    "D19" = "36.423.68"
    "E19" = "  -2.08%  "
    "D20" = "67.68"
    "E20" = "  -2.89%  "
    "E21" = "  -3.65%  "
    "D22" = "5.27"
    "E22" = "  +3.13%  "
    "D23" = "221.70"
    "E23" = "  -2.81%  "
    "E24" = "  +0.01%  "
    "E25" = "  +0.93%  "
    "E26" = "  -8.38%  "
    "D27" = "161.80"
    "E27" = "  -1.83%  "
    "D28" = "8.65"
    "E28" = "  -2.27%  "
    "E29" = "  -3.03%  "
    "D30" = "18.76"
    "E30" = "  -4.29%  "
    "D31" = "1.34"
    "E31" = "  +1.31%  "
    "D32" = "0.116"
    "E32" = "  -2.91%  "
    "D33" = "4.37"
    "E33" = "  -5.25%  "
    "D34" = "0.0601"
    "E34" = "  -6.52%  "
    "D35" = "4.24"
    "E35" = "  -5.92%  "
    "D36" = "2.33"
    "E36" = "  -1.10%  "
    "E37" = "  -0.01%  "
    "E38" = "  +0.65%  "
    "E39" = "  -2.84%  "
    "D40" = "5.61"
    "E40" = "  +6.60%  "
    "E41" = "  -1.41%  "
    "D42" = "0.0934"
    "E42" = "  +0.75%  "
    "D43" = "1.454.65"
    "E43" = "  +3.99%  "
    "E44" = "  -4.20%  "
    "D45" = "1.09"
    "E45" = "  -8.16%  "
    "D46" = "88.85"
    "E46" = "  -0.51%  "
    "D47" = "15.09"
    "E47" = "  -3.81%  "
    "D48" = "0.989"
    "E48" = "  -2.92%  "
    "D49" = "2.89"
    "E49" = "  -0.61%  "
    "D50" = "6.81"
    "E50" = "  -3.05%  "
    "D51" = "3.67"
    "E51" = "  +6.19%  "
}
foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
